{"js": "// Renumber every [[PERSON_N]] placeholder: for N >= 9, shift it down to\n// N-1 (PERSON_9 -> PERSON_8, PERSON_10 -> PERSON_9, ... PERSON_60 -> PERSON_59).\n// Placeholders PERSON_1 .. PERSON_8 are left untouched.\n//\n// We search the whole body for the wildcard pattern \"PERSON_[0-9]@\" (the\n// literal \"PERSON_\" followed by one or more digits) using Word's wildcard\n// search, which returns one range per occurrence of \"PERSON_<digits>\"\n// (the brackets \"[[ ]]\" around it are left alone). For each match we parse\n// out the number and, when it is >= 9, replace the whole matched range's\n// text with the decremented \"PERSON_<N-1>\" token.\n\nconst body = context.document.body;\nconst results = body.search(\"PERSON_[0-9]@\", { matchWildcards: true });\nresults.load(\"items/text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  const item = results.items[i];\n  const text = item.text;\n  const m = /^PERSON_(\\d+)$/.exec(text);\n  if (!m) {\n    continue;\n  }\n  const n = parseInt(m[1], 10);\n  if (n >= 9) {\n    item.insertText(\"PERSON_\" + (n - 1), Word.InsertLocation.replace);\n  }\n}\n\nawait context.sync();\n", "ps1": "# Renumber every [[PERSON_N]] placeholder: for N >= 9, shift it down to\n# N-1 (PERSON_9 -> PERSON_8, PERSON_10 -> PERSON_9, ... PERSON_60 -> PERSON_59).\n# Placeholders PERSON_1 .. PERSON_8 are left untouched.\n\n$d = $word.ActiveDocument\n\n# $rng is the Range that Find searches/mutates in place: after a successful\n# Execute() it is collapsed onto the just-found \"PERSON_<digits>\" text, so\n# reading/writing $rng.Text reads/replaces exactly that match.\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"PERSON_[0-9]@\"\n$find.MatchWildcards = $true\n$find.MatchCase = $true\n$find.Forward = $true\n$find.Wrap = 0\n\nwhile ($find.Execute()) {\n    $t = $rng.Text\n    if ($t -match '^PERSON_(\\d+)$') {\n        $n = [int]$matches[1]\n        if ($n -ge 9) {\n            $rng.Text = \"PERSON_\" + ($n - 1)\n        }\n    }\n}\n"}
